# =====================================================================
# Update the weekly CompStat report: new reporting period and refreshed
# crime-complaint statistics for precinct rows 15-30.
# =====================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: bump the volume/number and shift the reporting week ---
$ws.Range("A8").Value2 = "Volume 30   Number  42"
$ws.Range("C9").Value2 = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Donor cells used to clone the "no data" text style (style 14) ---
$zeroDonor = $ws.Cells.Item(14, 3)   # holds text "0"
$dashDonor = $ws.Cells.Item(14, 5)   # holds text "***.*"

# --- Row 15 ---
$c15 = $ws.Cells.Item(15, 3)
$c15.NumberFormat = '#,##0'
$c15.Value2 = 1
$d15 = $ws.Cells.Item(15, 4)
$d15.NumberFormat = '#,##0'
$d15.Value2 = 1
$e15 = $ws.Cells.Item(15, 5)
$e15.NumberFormat = '#,##0.0;"-"#,##0.0'
$e15.Value2 = 0
$g15 = $ws.Cells.Item(15, 7)
$g15.Value2 = 2
$h15 = $ws.Cells.Item(15, 8)
$h15.Value2 = 0
$i15 = $ws.Cells.Item(15, 9)
$i15.Value2 = 18
$j15 = $ws.Cells.Item(15, 10)
$j15.Value2 = 16
$k15 = $ws.Cells.Item(15, 11)
$k15.Value2 = 12.5
$l15 = $ws.Cells.Item(15, 12)
$l15.Value2 = 0
$m15 = $ws.Cells.Item(15, 13)
$m15.Value2 = 63.636363636363
$n15 = $ws.Cells.Item(15, 14)
$n15.Value2 = -43.75

# --- Row 16 ---
$c16 = $ws.Cells.Item(16, 3)
$c16.Value2 = 2
$d16 = $ws.Cells.Item(16, 4)
$d16.Value2 = 2
$f16 = $ws.Cells.Item(16, 6)
$f16.Value2 = 10
$h16 = $ws.Cells.Item(16, 8)
$h16.Value2 = -37.5
$i16 = $ws.Cells.Item(16, 9)
$i16.Value2 = 156
$j16 = $ws.Cells.Item(16, 10)
$j16.Value2 = 194
$k16 = $ws.Cells.Item(16, 11)
$k16.Value2 = -19.587628865979
$l16 = $ws.Cells.Item(16, 12)
$l16.Value2 = 12.230215827338
$m16 = $ws.Cells.Item(16, 13)
$m16.Value2 = -0.636942675159
$n16 = $ws.Cells.Item(16, 14)
$n16.Value2 = -79.922779922779

# --- Row 17 ---
$c17 = $ws.Cells.Item(17, 3)
$c17.Value2 = 4
$e17 = $ws.Cells.Item(17, 5)
$e17.Value2 = 0
$f17 = $ws.Cells.Item(17, 6)
$f17.Value2 = 12
$h17 = $ws.Cells.Item(17, 8)
$h17.Value2 = -29.411764705882
$i17 = $ws.Cells.Item(17, 9)
$i17.Value2 = 187
$j17 = $ws.Cells.Item(17, 10)
$j17.Value2 = 170
$k17 = $ws.Cells.Item(17, 11)
$k17.Value2 = 10
$l17 = $ws.Cells.Item(17, 12)
$l17.Value2 = 16.875
$m17 = $ws.Cells.Item(17, 13)
$m17.Value2 = 47.244094488189
$n17 = $ws.Cells.Item(17, 14)
$n17.Value2 = -60.042735042735

# --- Row 18 ---
$d18 = $ws.Cells.Item(18, 4)
$d18.Value2 = 3
$e18 = $ws.Cells.Item(18, 5)
$e18.Value2 = 33.333333333333
$f18 = $ws.Cells.Item(18, 6)
$f18.Value2 = 21
$g18 = $ws.Cells.Item(18, 7)
$g18.Value2 = 18
$h18 = $ws.Cells.Item(18, 8)
$h18.Value2 = 16.666666666666
$i18 = $ws.Cells.Item(18, 9)
$i18.Value2 = 212
$j18 = $ws.Cells.Item(18, 10)
$j18.Value2 = 318
$k18 = $ws.Cells.Item(18, 11)
$k18.Value2 = -33.333333333333
$l18 = $ws.Cells.Item(18, 12)
$l18.Value2 = -8.225108225108
$m18 = $ws.Cells.Item(18, 13)
$m18.Value2 = 17.127071823204
$n18 = $ws.Cells.Item(18, 14)
$n18.Value2 = -71.657754010695

# --- Row 19 ---
$c19 = $ws.Cells.Item(19, 3)
$c19.Value2 = 15
$e19 = $ws.Cells.Item(19, 5)
$e19.Value2 = -31.818181818181
$f19 = $ws.Cells.Item(19, 6)
$f19.Value2 = 75
$g19 = $ws.Cells.Item(19, 7)
$g19.Value2 = 89
$h19 = $ws.Cells.Item(19, 8)
$h19.Value2 = -15.730337078651
$i19 = $ws.Cells.Item(19, 9)
$i19.Value2 = 815
$j19 = $ws.Cells.Item(19, 10)
$j19.Value2 = 847
$k19 = $ws.Cells.Item(19, 11)
$k19.Value2 = -3.778040141676
$l19 = $ws.Cells.Item(19, 12)
$l19.Value2 = 33.606557377049
$m19 = $ws.Cells.Item(19, 13)
$m19.Value2 = 31.451612903225
$n19 = $ws.Cells.Item(19, 14)
$n19.Value2 = -37.595712098009

# --- Row 20 ---
$c20 = $ws.Cells.Item(20, 3)
$c20.Value2 = "0"
$zeroDonor.Copy($c20)
$d20 = $ws.Cells.Item(20, 4)
$d20.Value2 = 2
$e20 = $ws.Cells.Item(20, 5)
$e20.Value2 = -100
$f20 = $ws.Cells.Item(20, 6)
$f20.Value2 = 3
$g20 = $ws.Cells.Item(20, 7)
$g20.Value2 = 5
$h20 = $ws.Cells.Item(20, 8)
$h20.Value2 = -40
$j20 = $ws.Cells.Item(20, 10)
$j20.Value2 = 39
$k20 = $ws.Cells.Item(20, 11)
$k20.Value2 = -5.128205128205
$m20 = $ws.Cells.Item(20, 13)
$m20.Value2 = -9.756097560975
$n20 = $ws.Cells.Item(20, 14)
$n20.Value2 = -91.722595078299

# --- Row 21 ---
$d21 = $ws.Cells.Item(21, 4)
$d21.Value2 = 34
$e21 = $ws.Cells.Item(21, 5)
$e21.Value2 = -23.529411764705
$f21 = $ws.Cells.Item(21, 6)
$f21.Value2 = 123
$g21 = $ws.Cells.Item(21, 7)
$g21.Value2 = 147
$h21 = $ws.Cells.Item(21, 8)
$h21.Value2 = -16.326530612244
$i21 = $ws.Cells.Item(21, 9)
$i21.Value2 = 1426
$j21 = $ws.Cells.Item(21, 10)
$j21.Value2 = 1589
$k21 = $ws.Cells.Item(21, 11)
$k21.Value2 = -10.258023914411
$l21 = $ws.Cells.Item(21, 12)
$l21.Value2 = 14.538152610441
$m21 = $ws.Cells.Item(21, 13)
$m21.Value2 = 25.087719298245
$n21 = $ws.Cells.Item(21, 14)
$n21.Value2 = -62.384595093642

# --- Row 22 ---
$l22 = $ws.Cells.Item(22, 12)
$l22.Value2 = 12.5

# --- Row 23 ---
$c23 = $ws.Cells.Item(23, 3)
$c23.Value2 = 3
$d23 = $ws.Cells.Item(23, 4)
$d23.Value2 = 1
$e23 = $ws.Cells.Item(23, 5)
$e23.Value2 = 200
$f23 = $ws.Cells.Item(23, 6)
$f23.Value2 = 8
$h23 = $ws.Cells.Item(23, 8)
$h23.Value2 = 0
$i23 = $ws.Cells.Item(23, 9)
$i23.Value2 = 106
$j23 = $ws.Cells.Item(23, 10)
$j23.Value2 = 114
$k23 = $ws.Cells.Item(23, 11)
$k23.Value2 = -7.017543859649
$l23 = $ws.Cells.Item(23, 12)
$l23.Value2 = -36.526946107784
$m23 = $ws.Cells.Item(23, 13)
$m23.Value2 = 10.416666666666

# --- Row 24 ---
$c24 = $ws.Cells.Item(24, 3)
$c24.Value2 = 24
$d24 = $ws.Cells.Item(24, 4)
$d24.Value2 = 27
$e24 = $ws.Cells.Item(24, 5)
$e24.Value2 = -11.111111111111
$f24 = $ws.Cells.Item(24, 6)
$f24.Value2 = 129
$g24 = $ws.Cells.Item(24, 7)
$g24.Value2 = 140
$h24 = $ws.Cells.Item(24, 8)
$h24.Value2 = -7.857142857142
$i24 = $ws.Cells.Item(24, 9)
$i24.Value2 = 1220
$j24 = $ws.Cells.Item(24, 10)
$j24.Value2 = 1819
$k24 = $ws.Cells.Item(24, 11)
$k24.Value2 = -32.930181418361
$l24 = $ws.Cells.Item(24, 12)
$l24.Value2 = 33.333333333333
$m24 = $ws.Cells.Item(24, 13)
$m24.Value2 = -12.857142857142

# --- Row 25 ---
$f25 = $ws.Cells.Item(25, 6)
$f25.Value2 = 29
$g25 = $ws.Cells.Item(25, 7)
$g25.Value2 = 27
$h25 = $ws.Cells.Item(25, 8)
$h25.Value2 = 7.407407407407
$i25 = $ws.Cells.Item(25, 9)
$i25.Value2 = 375
$j25 = $ws.Cells.Item(25, 10)
$j25.Value2 = 387
$k25 = $ws.Cells.Item(25, 11)
$k25.Value2 = -3.100775193798
$l25 = $ws.Cells.Item(25, 12)
$l25.Value2 = 19.047619047619
$m25 = $ws.Cells.Item(25, 13)
$m25.Value2 = -2.088772845953

# --- Row 26 ---
$c26 = $ws.Cells.Item(26, 3)
$c26.NumberFormat = '#,##0'
$c26.Value2 = 2
$d26 = $ws.Cells.Item(26, 4)
$d26.NumberFormat = '#,##0'
$d26.Value2 = 2
$e26 = $ws.Cells.Item(26, 5)
$e26.NumberFormat = '#,##0.0;"-"#,##0.0'
$e26.Value2 = 0
$g26 = $ws.Cells.Item(26, 7)
$g26.Value2 = 5
$h26 = $ws.Cells.Item(26, 8)
$h26.Value2 = -40
$i26 = $ws.Cells.Item(26, 9)
$i26.Value2 = 31
$j26 = $ws.Cells.Item(26, 10)
$j26.Value2 = 30
$k26 = $ws.Cells.Item(26, 11)
$k26.Value2 = 3.333333333333
$l26 = $ws.Cells.Item(26, 12)
$l26.Value2 = 14.814814814814

# --- Row 27 ---
$f27 = $ws.Cells.Item(27, 6)
$f27.Value2 = 4
$j27 = $ws.Cells.Item(27, 10)
$j27.Value2 = 79
$k27 = $ws.Cells.Item(27, 11)
$k27.Value2 = -41.772151898734

# --- Row 28 ---
$d28 = $ws.Cells.Item(28, 4)
$d28.NumberFormat = '#,##0'
$d28.Value2 = 1
$e28 = $ws.Cells.Item(28, 5)
$e28.NumberFormat = '#,##0.0;"-"#,##0.0'
$e28.Value2 = -100
$g28 = $ws.Cells.Item(28, 7)
$g28.NumberFormat = '#,##0'
$g28.Value2 = 1
$h28 = $ws.Cells.Item(28, 8)
$h28.NumberFormat = '#,##0.0;"-"#,##0.0'
$h28.Value2 = -100
$j28 = $ws.Cells.Item(28, 10)
$j28.Value2 = 6
$k28 = $ws.Cells.Item(28, 11)
$k28.Value2 = -50
$n28 = $ws.Cells.Item(28, 14)
$n28.Value2 = -89.285714285714

# --- Row 29 ---
$d29 = $ws.Cells.Item(29, 4)
$d29.NumberFormat = '#,##0'
$d29.Value2 = 1
$e29 = $ws.Cells.Item(29, 5)
$e29.NumberFormat = '#,##0.0;"-"#,##0.0'
$e29.Value2 = -100
$g29 = $ws.Cells.Item(29, 7)
$g29.NumberFormat = '#,##0'
$g29.Value2 = 1
$h29 = $ws.Cells.Item(29, 8)
$h29.NumberFormat = '#,##0.0;"-"#,##0.0'
$h29.Value2 = -100
$j29 = $ws.Cells.Item(29, 10)
$j29.Value2 = 6
$k29 = $ws.Cells.Item(29, 11)
$k29.Value2 = -50
$n29 = $ws.Cells.Item(29, 14)
$n29.Value2 = -86.95652173913

# --- Row 30 ---
$g30 = $ws.Cells.Item(30, 7)
$g30.Value2 = "0"
$zeroDonor.Copy($g30)
$h30 = $ws.Cells.Item(30, 8)
$h30.Value2 = "***.*"
$dashDonor.Copy($h30)
